# demonstrating setting inputs from file
#
# Adds a "jurisdiction" lookup column to the travel and relative-mixing
# matrices (so each row can pull its row label from a single shared value),
# and extends the parameters sheet with a "source" column and an extra
# parameter row ("…") to show settings can be sourced from file. Finally
# moves the active tab to the "parameters" sheet.

$wb = $excel.ActiveWorkbook

# --- travel sheet: add jurisdiction label in A1, reset selection ---
$wsTravel = $wb.Worksheets.Item("travel")
$wsTravel.Range("A1").Value = "jurisdiction"
[void]$wsTravel.Range("A2").Select()

# --- relative-mixing sheet: add jurisdiction label in A1, reset selection ---
$wsMixing = $wb.Worksheets.Item("relative-mixing")
$wsMixing.Range("A1").Value = "jurisdiction"
[void]$wsMixing.Range("A2").Select()

# --- parameters sheet: add a new "…" row and "source" column header ---
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("A4").Value = [char]0x2026
$wsParams.Range("F1").Value = "source"

# Make "parameters" the active sheet/tab and select the full used range.
[void]$wsParams.Activate()
[void]$wsParams.Range("A1:F4").Select()
